# Sales/ER-wait-time QI workbook: the baseline/center-line calculation was
# switched from the mean to the median, and the chart-support label column
# was relabeled to match ("ave" -> "median"; the pre-existing "baseline"
# label in C1 is unaffected content-wise, it just keeps pointing at the same
# shared string). D (UCL) / E (LCL) recompute automatically since they are
# formulas driven off H3/H2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("L1")

# G3 used to read "ave" (paired with H3 = AVERAGE(B:B)); rename to "median"
# and switch the statistic itself to MEDIAN.
$ws.Range("G3").Value = "median"
$ws.Range("H3").Formula = "=MEDIAN(B:B)"

# Force a full recalculation so D:E (UCL/LCL, which depend on H2:H3) and
# every other dependent formula cell picks up the new median-based value.
$wb.Application.Calculate()

# The author had last clicked on H10 before saving.
$ws.Range("H10").Select()
